$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.993.62'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.872.61'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'305.64"
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = "'0.5082"
$ws.Range('D8').Value = "'0.3663"
$ws.Range('E8').Value = '  -2.41%  '
$ws.Range('D9').Value = "'0.07217"
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = "'0.8955"
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = "'95.32"
$ws.Range('E13').Value = '  +6.53%  '
$ws.Range('D14').Value = '1.872.15'
$ws.Range('E14').Value = '  +0.37%  '
$ws.Range('D15').Value = "'5.249"
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D17').Value = "'0.000008539"
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = "'14.25"
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '27.018.91'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('D22').Value = '2.102.01'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').Value = "'10.39"
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').Value = "'6.400"
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').Value = "'148.60"
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('D27').Value = "'17.91"
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').Value = "'2.093"
$ws.Range('E28').Value = '  -1.10%  '
$ws.Range('D29').Value = "'113.38"
$ws.Range('E29').Value = '  +0.47%  '
$ws.Range('D30').Value = "'4.724"
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('D31').Value = "'4.750"
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').Value = "'0.09182"
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').Value = "'0.05108"
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').Value = "'0.7515"
$ws.Range('E34').Value = '  +3.24%  '
$ws.Range('D35').Value = "'2.976"
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').Value = "'1.159"
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  +6.29%  '
$ws.Range('D38').Value = "'2.551"
$ws.Range('E38').Value = '  +3.03%  '
$ws.Range('D39').Value = "'0.5656"
$ws.Range('E39').Value = '  +6.04%  '
$ws.Range('D40').Value = "'0.02001"
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('D41').Value = "'1.077"
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = "'6.643"
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').Value = "'115.70"
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').Value = "'8.574"
$ws.Range('E44').Value = '  +3.59%  '
$ws.Range('D45').Value = "'0.1478"
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('D46').Value = "'0.4771"
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = "'1.000"
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'10.12"
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = "'36.97"
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = "'63.14"
$ws.Range('E51').Value = '  -1.26%  '
